$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Refraction os Solid")

for ($i = 1; $i -le 15; $i++) {
    $ws.Cells.Item($i, 3).Value = $i
}

$ws.Activate()
$ws.Range("C1:C15").Select()
